$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.321.32"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "1.856.49"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("E4").Value = "  -0.76%  "
$ws.Range("D5").Value = "'314.18"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("D7").Value = "'0.4611"
$ws.Range("E7").Value = "  -0.98%  "
$ws.Range("D8").Value = "'0.3702"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("D9").Value = "'0.07338"
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("D10").Value = "'0.8830"
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("D11").Value = "'19.87"
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("D12").Value = "'0.07804"
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("D13").Value = "1.876.63"
$ws.Range("E13").Value = "  +2.29%  "
$ws.Range("D14").Value = "'5.390"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").Value = "'6.547"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").Value = "'91.89"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").Value = "'0.000009014"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("D20").Value = "'14.79"
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("D21").Value = "27.342.60"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("D22").Value = "'5.123"
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").Value = "'10.48"
$ws.Range("E23").Value = "  -1.24%  "
$ws.Range("D24").Value = "2.096.22"
$ws.Range("E24").Value = "  +1.60%  "
$ws.Range("D25").Value = "'1.915"
$ws.Range("E25").Value = "  +4.22%  "
$ws.Range("D26").Value = "'152.05"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").Value = "'18.36"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").Value = "'2.072"
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("D29").Value = "'5.118"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").Value = "'116.13"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").Value = "'0.08850"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").Value = "'0.7670"
$ws.Range("E32").Value = "  +5.07%  "
$ws.Range("D33").Value = "'3.002"
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("D34").Value = "'1.175"
$ws.Range("E34").Value = "  +3.40%  "
$ws.Range("D35").Value = "'4.497"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("D36").Value = "'2.639"
$ws.Range("E36").Value = "  +6.26%  "
$ws.Range("D37").Value = "'0.01962"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").Value = "'1.078"
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("D39").Value = "'0.05228"
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("D40").Value = "'2.943"
$ws.Range("E40").Value = "  +0.51%  "
$ws.Range("D41").Value = "'7.051"
$ws.Range("E41").Value = "  -4.48%  "
$ws.Range("D42").Value = "'0.5148"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").Value = "'0.1639"
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("D44").Value = "'8.365"
$ws.Range("E44").Value = "  +1.47%  "
$ws.Range("D45").Value = "'0.4837"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("D46").Value = "'10.30"
$ws.Range("E46").Value = "  +0.82%  "
$ws.Range("D47").Value = "'1.000"
$ws.Range("D48").Value = "'103.10"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("D49").Value = "'1.654"
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("D50").Value = "'0.06218"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("D51").Value = "'65.82"
$ws.Range("E51").Value = "  +1.99%  "
